# Update the "actual" progress markers (who has finished which task) and
# mark the PowerPoint/presentation task as completed with an actual end date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C7: Window Form C# research is now done by "Tran Lieu Nhut Anh" -> mark X
$ws.Range("C7").Value = "X"

# C10 / C11: these two items are no longer marked done for that person
$ws.Range("C10").Value = ""
$ws.Range("C11").Value = ""

# D13 / D14 / D15: these items are no longer marked done for that person
$ws.Range("D13").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("D15").Value = ""

# C16: no longer marked done
$ws.Range("C16").Value = ""

# H17: presentation task now has an actual completion date instead of "Chưa kết thúc"
$ws.Range("H17").Value = 43779

# Move the active cell selection to I13 (matches the saved view state)
$ws.Range("I13").Select()
